$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'schubert-winterreise_195'
$ws.Range("B2").Value = 'isophonics_138'
$ws.Range("C2").Value = 0.111969111969112
$ws.Range("D2").Value = '[[''C:7'', ''F:min'', ''A#:min/C#'']]'
$ws.Range("E2").Value = '[[''G:7'', ''C:min'', ''F:min'']]'
$ws.Range("F2").Value = '[(23.4, 31.32)]'
$ws.Range("G2").Value = '[(2.352018, 8.679455)]'
$ws.Range("I2").Value = 'spotify:track:6tQvjqDIK9GXWIC6mejms8'

# Row 3
$ws.Range("A3").Value = 'jaah_16'
$ws.Range("B3").Value = 'jaah_39'
$ws.Range("C3").Value = 0.06742294520547945
$ws.Range("D3").Value = '[[''B:7'', ''E:7'', ''A'']]'
$ws.Range("E3").Value = '[[''Bb:7'', ''Eb:7'', ''Ab'']]'
$ws.Range("F3").Value = '[(42.96, 45.78)]'
$ws.Range("G3").Value = '[(48.18, 50.27)]'

# Row 4
$ws.Range("A4").Value = 'jaah_44'
$ws.Range("B4").Value = 'jaah_0'
$ws.Range("C4").Value = 0.1017241379310345
$ws.Range("D4").Value = '[[''Ab:7'', ''Db'', ''Db'']]'
$ws.Range("E4").Value = '[[''Eb:7'', ''Ab'', ''Ab'']]'
$ws.Range("F4").Value = '[(46.99, 49.5)]'
$ws.Range("G4").Value = '[(64.32, 68.27)]'

# Row 5
$ws.Range("A5").Value = 'isophonics_76'
$ws.Range("B5").Value = 'isophonics_223'
$ws.Range("C5").Value = 0.1112637362637363
$ws.Range("D5").Value = '[[''E:min'', ''G:7'', ''C''], [''G'', ''D/3'', ''E:min'']]'
$ws.Range("E5").Value = '[[''F#:min'', ''A:7'', ''D''], [''A'', ''E'', ''F#:min'']]'
$ws.Range("F5").Value = '[(3.733, 8.605), (10.141, 14.901)]'
$ws.Range("G5").Value = '[(11.929818, 15.436031), (9.607823, 12.823786)]'
$ws.Range("H5").Value = ''
$ws.Range("I5").Value = 'spotify:track:3KfbEIOC7YIv90FIfNSZpo'

# Row 6
$ws.Range("A6").Value = 'isophonics_292'
$ws.Range("B6").Value = 'isophonics_96'
$ws.Range("C6").Value = 0.2864583333333334
$ws.Range("D6").Value = '[[''E'', ''D/2'', ''A/5'', ''E'', ''D'']]'
$ws.Range("E6").Value = '[[''D:maj'', ''C:maj'', ''G/3'', ''D:maj'', ''C:maj'']]'
$ws.Range("F6").Value = '[(1.701428, 7.007188)]'
$ws.Range("G6").Value = '[(73.89, 86.025)]'

# Row 7
$ws.Range("A7").Value = 'jaah_54'
$ws.Range("B7").Value = 'jaah_85'
$ws.Range("C7").Value = 0.119106699751861
$ws.Range("D7").Value = '[[''Eb'', ''Eb'', ''Bb:7'', ''Bb:7'', ''Eb'']]'
$ws.Range("E7").Value = '[[''Ab'', ''Ab'', ''Eb:7'', ''Eb:7'', ''Ab'']]'
$ws.Range("F7").Value = '[(22.37, 40.37)]'
$ws.Range("G7").Value = '[(85.75, 90.9)]'
$ws.Range("H7").Value = ''

# Row 8
$ws.Range("A8").Value = 'schubert-winterreise_93'
$ws.Range("B8").Value = 'schubert-winterreise_178'
$ws.Range("C8").Value = 0.150268336314848
$ws.Range("D8").Value = '[[''C#:maj/G#'', ''G#:7'', ''C#:maj'']]'
$ws.Range("E8").Value = '[[''F:maj'', ''C:7/E'', ''F:maj'']]'
$ws.Range("F8").Value = '[(36.6, 39.0)]'
$ws.Range("G8").Value = '[(39.02, 43.96)]'
$ws.Range("H8").Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'

# Row 9
$ws.Range("A9").Value = 'isophonics_28'
$ws.Range("B9").Value = 'isophonics_196'
$ws.Range("C9").Value = 0.259958071278826
$ws.Range("D9").Value = '[[''G'', ''C/5'', ''G'', ''C''], [''C'', ''G'', ''C'', ''C'']]'
$ws.Range("E9").Value = '[[''A'', ''D/5'', ''A'', ''D/5''], [''D/5'', ''A'', ''D/5'', ''D'']]'
$ws.Range("F9").Value = '[(64.476575, 70.099405), (1.620158, 15.8656)]'
$ws.Range("G9").Value = '[(31.23721, 34.035215), (34.743424, 38.330907)]'
$ws.Range("H9").Value = ''

# Row 10
$ws.Range("A10").Value = 'isophonics_64'
$ws.Range("B10").Value = 'isophonics_128'
$ws.Range("C10").Value = 0.2302631578947368
$ws.Range("D10").Value = '[[''A'', ''D'', ''A'', ''D'']]'
$ws.Range("E10").Value = '[[''C'', ''F'', ''C'', ''F'']]'
$ws.Range("F10").Value = '[(122.78, 137.555)]'
$ws.Range("G10").Value = '[(10.634761, 19.597664)]'
$ws.Range("H10").Value = ''

# Row 11
$ws.Range("A11").Value = 'schubert-winterreise_161'
$ws.Range("B11").Value = 'isophonics_297'
$ws.Range("C11").Value = 0.3059210526315789
$ws.Range("D11").Value = '[[''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'']]'
$ws.Range("E11").Value = '[[''G'', ''D'', ''G'', ''D'', ''G'', ''D'', ''G'', ''D'', ''G'']]'
$ws.Range("F11").Value = '[(11.86, 17.78)]'
$ws.Range("G11").Value = '[(0.421247, 8.635573)]'

# Row 12
$ws.Range("A12").Value = 'schubert-winterreise_192'
$ws.Range("B12").Value = 'schubert-winterreise_66'
$ws.Range("C12").Value = 0.09642857142857142
$ws.Range("D12").Value = '[[''F:min/C'', ''C'', ''F:min/C'']]'
$ws.Range("E12").Value = '[[''C:min/G'', ''G'', ''C:min'']]'
$ws.Range("F12").Value = '[(47.68, 51.98)]'
$ws.Range("G12").Value = '[(12.28, 17.3)]'
$ws.Range("I12").Value = ''

# Row 13
$ws.Range("A13").Value = 'schubert-winterreise_44'
$ws.Range("B13").Value = 'schubert-winterreise_51'
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = '[[''A#:min'', ''A#:min'', ''C:hdim7/A#'', ''A#:min'', ''A#:min/F'', ''F:7'', ''A#:min'', ''A#:min'', ''C:hdim7/A#'', ''A#:min'', ''A#:min/F'', ''F:7'', ''A#:min'', ''G#/C'', ''G#:7'', ''C#'', ''A#:min'', ''D#:7/A#'', ''G#'', ''C#/F'', ''G#:7/D#'', ''C#/F'', ''G#:7'', ''C#'', ''C:dim7/C#'', ''C#:7'', ''F#'', ''D#:min'', ''G#:7/D#'', ''C#'', ''F#/A#'', ''C#:7/G#'', ''F#/A#'', ''C#'', ''F#'', ''C:hdim7/D#'', ''A#:min/C#'', ''C:hdim7/D#'', ''A#:min/C#'', ''C#/F'', ''A#:min/F'', ''F:7'', ''A#:min'', ''C:hdim7/D#'', ''A#:min/C#'', ''C#/F'', ''A#:min/F'', ''F:7'', ''A#:min'', ''A:dim7/A#'', ''A#:min'', ''D#:min'', ''A#:min/F'', ''D#:min7/F#'', ''A#:min/F'', ''F:7'', ''A#:min'', ''A#:min'', ''C:hdim7/A#'', ''A#:min'']]'
$ws.Range("E13").Value = '[[''C:min'', ''C:min'', ''D:hdim7/C'', ''C:min'', ''C:min/G'', ''G:7'', ''C:min'', ''C:min'', ''D:hdim7/C'', ''C:min'', ''C:min/G'', ''G:7'', ''C:min'', ''A#/D'', ''A#:7'', ''D#'', ''C:min'', ''F:7/C'', ''A#'', ''D#/G'', ''A#:7/F'', ''D#/G'', ''A#:7'', ''D#'', ''D:dim7/D#'', ''D#:7'', ''G#'', ''F:min'', ''A#:7/F'', ''D#'', ''G#/C'', ''D#:7/A#'', ''G#/C'', ''D#'', ''G#'', ''D:hdim7/F'', ''C:min/D#'', ''D:hdim7/F'', ''C:min/D#'', ''D#/G'', ''C:min/G'', ''G:7'', ''C:min'', ''D:hdim7/F'', ''C:min/D#'', ''D#/G'', ''C:min/G'', ''G:7'', ''C:min'', ''B:dim7/C'', ''C:min'', ''F:min'', ''C:min/G'', ''F:min7/G#'', ''C:min/G'', ''G:7'', ''C:min'', ''C:min'', ''D:hdim7/C'', ''C:min'']]'
$ws.Range("F13").Value = '[(15.9, 106.38)]'
$ws.Range("G13").Value = '[(14.88, 99.5)]'
$ws.Range("I13").Value = ''

# Row 14
$ws.Range("A14").Value = 'isophonics_165'
$ws.Range("B14").Value = 'schubert-winterreise_96'
$ws.Range("C14").Value = 0.115546218487395
$ws.Range("D14").Value = '[[''A'', ''E'', ''A'']]'
$ws.Range("E14").Value = '[[''F:maj'', ''C:maj'', ''F:maj'']]'
$ws.Range("F14").Value = '[(52.680839, 55.420793)]'
$ws.Range("G14").Value = '[(0.74, 21.5)]'

# Row 15
$ws.Range("A15").Value = 'isophonics_109'
$ws.Range("B15").Value = 'schubert-winterreise_205'
$ws.Range("C15").Value = 0.5230769230769231
$ws.Range("D15").Value = '[[''A'', ''D'', ''A'']]'
$ws.Range("E15").Value = '[[''G:maj'', ''C:maj/G'', ''G:maj'']]'
$ws.Range("F15").Value = '[(94.925, 100.821)]'
$ws.Range("G15").Value = '[(16.92, 23.9)]'
$ws.Range("I15").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_44'
$ws.Range("B16").Value = 'schubert-winterreise_88'
$ws.Range("C16").Value = 0.1613636363636363
$ws.Range("D16").Value = '[[''A#/F'', ''F:7'', ''A#''], [''A#'', ''A#/F'', ''F:7'']]'
$ws.Range("E16").Value = '[[''A:maj/E'', ''E:7'', ''A:maj''], [''A:maj'', ''A:maj'', ''E:7/G#'']]'
$ws.Range("F16").Value = '[(271.22, 275.78), (269.7, 272.54)]'
$ws.Range("G16").Value = '[(16.04, 18.72), (6.56, 12.6)]'
$ws.Range("I16").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

# Row 17
$ws.Range("A17").Value = 'schubert-winterreise_6'
$ws.Range("B17").Value = 'schubert-winterreise_180'
$ws.Range("C17").Value = 0.09018567639257294
$ws.Range("D17").Value = '[[''B:min/F#'', ''F#:7'', ''B:min''], [''B:min'', ''F#:maj'', ''B:min'']]'
$ws.Range("E17").Value = '[[''G:min'', ''D:7/G'', ''G:min''], [''G:min'', ''D:maj/G'', ''G:min'']]'
$ws.Range("F17").Value = '[(79.42, 86.02), (13.98, 21.5)]'
$ws.Range("G17").Value = '[(12.06, 17.06), (4.84, 10.64)]'
$ws.Range("H17").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

Write-Output "applied edits"
